# "edit export excel and fixing status column => still cannot render properly"
#
# The Daily_History sheet previously held three rows (2-4) of placeholder
# multiplication-table numbers (2,4,6,8... / 3,6,9,12... / 4,8,12,16...).
# This replaces that placeholder data with a single real sample row and
# removes the now-unused rows 3 and 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily_History")

# Drop the two extra placeholder rows - only one sample data row remains.
$ws.Rows("3:4").Delete()

# Clear out the old placeholder numbers in row 2 before writing new data.
$ws.Range("A2:P2").ClearContents()

# Columns that hold numeric-looking text (dates / ids stored as strings,
# not numbers) need to be forced to Text format first so Excel keeps them
# as strings instead of silently coercing them to numbers.
$textCells = @("A2", "B2", "D2", "E2", "M2", "O2")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("A2").Value = "20200313"
$ws.Range("B2").Value = "999999"
$ws.Range("C2").Value = "Other"
$ws.Range("D2").Value = "20150531"
$ws.Range("E2").Value = "20160531"
$ws.Range("F2").Value = "undefined"
$ws.Range("G2").Value = 180
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 30
$ws.Range("K2").Value = 13
$ws.Range("L2").Value = 30
$ws.Range("M2").Value = "1"
$ws.Range("N2").Value = "Akiyama"
$ws.Range("O2").Value = "41"
$ws.Range("P2").Value = "Training/Support"

# Now that the text values are committed, drop the temporary Text number
# format again so the cells end up back at the workbook's default style.
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
